# DESCW-1269 Change Request Types (Tab 22)
# Build out the report template: replace the placeholder title with the
# real report title, and fill in the header/subheader/detail-row template
# tokens for the Change-Request-Types table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (row 1, merged C1:I1) ---------------------------------------
$ws.Range("C1").Value = 'Change Request Types for {$fy} as of {$date}'

# --- Column headers (row 2) ---------------------------------------------
$ws.Range("C2").Value = 'FY'
$ws.Range("D2").Value = 'Total CRs'
$ws.Range("E2").Value = 'Initiated By'
$ws.Range("F2").Value = 'Budget'
$ws.Range("G2").Value = 'Schedule'
$ws.Range("H2").Value = 'Scope'
$ws.Range("I2").Value = 'None'

# --- Detail-row template tokens (row 3) ----------------------------------
$ws.Range("A3").Value = '{$r.project_number}'
$ws.Range("B3").Value = '{$r.project_name}'
$ws.Range("C3").Value = '{$r.fiscal_year}'
$ws.Range("D3").Value = '{$r.cr_count}'
$ws.Range("E3").Value = '{$r.initiated_by}'
$ws.Range("F3").Value = '{$r.budget}'
$ws.Range("G3").Value = '{$r.schedule}'
$ws.Range("H3").Value = '{$r.scope}'
$ws.Range("I3").Value = '{$r.none}'

# --- Second report-row loop token (row 4) --------------------------------
$ws.Range("A4").Value = '{$r1}'

# --- Move the active selection to B10 (was B7:B12) -----------------------
$ws.Range("B10").Select()
